$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("design_variables")

$ws.Range("C1").Value = 697674.41860465112
$ws.Range("B2").Value = 1932.558139534884
